# Project brief: annotate the "Hosting & domain name" budget line with
# "(for one year)" — it was previously followed immediately by a run of
# tab characters (used to push the "$   400" figure over to the right).
# The edit collapses the first two of those tab characters into the new
# text " (for one year)", leaving the remaining tabs (and the "$   400"
# amount) untouched.

$d = $word.ActiveDocument

# Build a two-tab-character string (PowerShell won't concatenate [char]
# values with '+' the way you might expect - that does numeric/char
# arithmetic - so use string interpolation instead).
$twoTabs = "$([char]9)$([char]9)"

# Locate the unique anchor text "Hosting & domain name" immediately
# followed by two tabs, searched across the whole document body.
$anchor = $d.Content
$found = $anchor.Find.Execute("Hosting & domain name" + $twoTabs, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # $anchor now spans the matched text ("Hosting & domain name" + the
    # two tabs). Narrow it down to just the trailing two tab characters
    # so the rest of the run/paragraph is left completely alone.
    $tabsStart = $anchor.End - 2
    $tabsEnd = $anchor.End
    $toReplace = $d.Range($tabsStart, $tabsEnd)
    $toReplace.Text = " (for one year)"
}
